$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.114.00"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").Value = "3.796.18"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'602.83"
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("D6").Value = "'165.71"
$ws.Range("E6").Value = "  -0.75%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.518"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("E11").Value = "  +1.50%  "

$ws.Range("E12").Value = "  -1.21%  "

$ws.Range("D13").Value = "'35.93"
$ws.Range("E13").Value = "  -0.24%  "

$ws.Range("D14").Value = "4.431.88"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").Value = "3.792.33"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").Value = "68.102.03"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").Value = "'18.31"
$ws.Range("E17").Value = "  -1.48%  "

$ws.Range("E18").Value = "  +1.95%  "

$ws.Range("D19").Value = "'7.09"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "'462.70"
$ws.Range("E20").Value = "  +0.41%  "

$ws.Range("D21").Value = "'9.73"
$ws.Range("E21").Value = "  -2.19%  "

$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("D23").Value = "'0.0000149"
$ws.Range("E23").Value = "  -2.45%  "

$ws.Range("D24").Value = "'82.94"
$ws.Range("E24").Value = "  -0.50%  "

$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").Value = "'2.11"
$ws.Range("E26").Value = "  +0.69%  "

$ws.Range("D28").Value = "'10.00"
$ws.Range("E28").Value = "  +0.20%  "

$ws.Range("D29").Value = "3.944.87"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "'7.45"
$ws.Range("E30").Value = "  +3.04%  "

$ws.Range("D31").Value = "'2.64"
$ws.Range("E31").Value = "  -4.89%  "

$ws.Range("E32").Value = "  -1.55%  "

$ws.Range("D33").Value = "'29.31"
$ws.Range("E33").Value = "  -1.05%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").Value = "'9.01"
$ws.Range("E35").Value = "  -0.44%  "

$ws.Range("D36").Value = "'0.0997"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  +1.18%  "

$ws.Range("D38").Value = "'3.28"
$ws.Range("E38").Value = "  -2.24%  "

$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").Value = "'0.987"
$ws.Range("E40").Value = "  -0.72%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.02%  "

$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("D44").Value = "'47.45"
$ws.Range("E44").Value = "  -1.47%  "

$ws.Range("D45").Value = "'43.18"
$ws.Range("E45").Value = "  -1.53%  "

$ws.Range("D46").Value = "'150.98"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("E47").Value = "  +0.77%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.86"
$ws.Range("E48").Value = "  +2.67%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'395.61"
$ws.Range("E49").Value = "  +1.21%  "

$ws.Range("D50").Value = "'26.95"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("D51").Value = "'1.33"
$ws.Range("E51").Value = "  +5.26%  "
